# JournalDeBord.xlsx edit:
#  - Insert a new log row after the "Création d'images bitmap..." entry (old row 29)
#    for a new "Lecture de documentation sur les événements tkinter" entry (0.5h).
#  - Update the wording of the "Création d'images bitmap pour l'interface" entry.
#  - Fill in the hours (1h) for that entry.
#  - Extend the print area / dimension accordingly (handled automatically by Excel
#    once the new row with data exists, but we set it explicitly to be safe).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# 1. Update wording of the existing row 29 entry and fill in its hours.
$ws.Range("B29").Value = "Création d'images bitmap pour l'interface et intégration dans le code"
$ws.Range("C29").Value = 1

# 2. Insert a new row before row 30 (shifts old rows 30-41 down to 31-42).
$ws.Rows.Item(30).Insert(-4121, 0)   # xlShiftDown = -4121, xlFormatFromLeftOrAbove = 0

# 2b. Re-apply the formatting of the row above (row 29) onto the freshly
#     inserted, still-blank row so its cell styles (borders, number formats)
#     match the rest of the table exactly, same as Excel's own behaviour
#     when inserting a row in the middle of a formatted table.
$ws.Range("A29:D29").Copy()
$ws.Range("A30:D30").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# 3. Populate the newly inserted row 30 with the new journal entry.
$ws.Range("B30").Value = "Lecture de documentation sur les événements tkinter"
$ws.Range("C30").Value = 0.5

# 4. Extend the print area by one row to include the new total row.
$ws.PageSetup.PrintArea = '$A$1:$C$42'

# 5. Move the active selection, matching the post-edit workbook state.
$ws.Range("H24").Select()

$wb.Save()
